$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10, pushing the old rows 10 and 11 down to 11 and 12
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new weekly data
$ws.Range("A10").Value = 12
$ws.Range("B10").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C10").Value = "Metropolitana"
$ws.Range("D10").Value = 44460
$ws.Range("E10").Value = 13
$ws.Range("F10").Value = 100112013
$ws.Range("G10").Value = "Alcachofa"
$ws.Range("H10").Value = "Española"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 45
$ws.Range("K10").Value = 13000
$ws.Range("L10").Value = 13000
$ws.Range("M10").Value = 13000
$ws.Range("N10").Value = '$/caja 30 unidades'
$ws.Range("O10").Value = "Provincia de Limarí"
$ws.Range("P10").Value = 433
$ws.Range("Q10").Value = 30
$ws.Range("R10").Value = "Hortaliza"
